$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values (coordinate-based writes) ---
$ws.Cells.Item(1, 1).Value = "Testing"

$ws.Cells.Item(2, 1).Value = "A"
$ws.Cells.Item(2, 2).Value = "B"
$ws.Cells.Item(2, 3).Value = "C"

$ws.Cells.Item(3, 1).Value = "X"
$ws.Cells.Item(3, 7).Value = "CHI"

$ws.Cells.Item(4, 1).Value = "Y"
$ws.Cells.Item(4, 3).Value = "ALPHA"
$ws.Cells.Item(4, 4).Value = "BETA"
$ws.Cells.Item(4, 5).Value = "GAMMA"
$ws.Cells.Item(4, 7).Value = "PSI"

$ws.Cells.Item(5, 1).Value = "Z"
$ws.Cells.Item(5, 7).Value = "OMEGA"

# --- Header-style formatting (bold, centered, grey fill) for the new label cells ---
# Reuse A1's existing header style (bold / centered / grey fill) via copy-format
# so the cells share the same style record instead of minting new ones.
$ws.Range("A1").Copy()
$headerCells = @("G3", "C4", "D4", "E4", "G4", "G5")
foreach ($addr in $headerCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.9
$ws.Columns.Item(2).ColumnWidth = 25

# --- Freeze panes at B2 (1 column / 1 row frozen) ---
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Page margins (header/footer = 0.5in) ---
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
